$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 3000
$ws.Range("B2").Value = 3000
$ws.Range("C2").Value = 20

$ws.Range("A3").Value = 4000
$ws.Range("B3").Value = 4000
$ws.Range("C3").Value = 20

$ws.Range("A4").Value = 1000
$ws.Range("B4").Value = 1000
$ws.Range("C4").Value = 10

$ws.Range("A5").Value = 4000
$ws.Range("B5").Value = 1000
$ws.Range("C5").Value = 10
